$d = $word.ActiveDocument

$replacements = @(
    @("2025-10-30 Thursday", "2025-10-31 Friday"),
    @("801×2=", "846×4="),
    @("458×6=", "172×9="),
    @("559×3=", "398×4="),
    @("888×4=", "419×7="),
    @("567×5=", "613×7="),
    @("191×2=", "965×4="),
    @("112×3=", "393×5="),
    @("252×6=", "179×9="),
    @("377×9=", "941×6="),
    @("138×3=", "243×2="),
    @("227×3=", "494×4="),
    @("520×5=", "750×8="),
    @("416×3=", "749×4="),
    @("279×3=", "621×5="),
    @("835×7=", "750×2="),
    @("480×4=", "416×2="),
    @("233×8=", "242×4="),
    @("485×7=", "580×2="),
    @("261×5=", "584×2="),
    @("617×3=", "782×6="),
    @("672×4=", "816×4="),
    @("431×6=", "131×2="),
    @("226×6=", "748×6="),
    @("410×5=", "709×4="),
    @("277×7=", "787×7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
